$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the header row
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 05:45"

# Belgica (row 40) - updated counts, no reordering needed
$ws.Range("B40").Value = 85911
$ws.Range("C40").Value = 424
$ws.Range("D40").Value = 18490
$ws.Range("E40").Value = 57523
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 9898

# Lesoto / Belice swap places (rows 162-163): Belice overtakes Lesoto
$ws.Range("A162").Value = "Belice"
$ws.Range("B162").Value = 1101
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 240
$ws.Range("E162").Value = 848
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 13

$ws.Range("A163").Value = "Lesoto"
$ws.Range("B163").Value = 1085
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 528
$ws.Range("E163").Value = 526
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 31

# Taiwan / San Martin (Parte Holandesa) swap places (rows 174-175)
$ws.Range("A174").Value = "San Martin (Parte Holandesa)"
$ws.Range("B174").Value = 495
$ws.Range("C174").Value = 13
$ws.Range("D174").Value = 302
$ws.Range("E174").Value = 174
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 19

$ws.Range("A175").Value = "Taiwan"
$ws.Range("B175").Value = 489
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 471
$ws.Range("E175").Value = 11
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 7

# Camboya (row 185) - updated counts, no reordering needed
$ws.Range("D185").Value = 271
$ws.Range("E185").Value = 3

# Curazao (row 196) - updated counts, no reordering needed
$ws.Range("B196").Value = 75
$ws.Range("C196").Value = 4
$ws.Range("D196").Value = 39

# Montserrat / Islas Malvinas swap places (rows 214-215)
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
